$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case row marker in column A of row 2 (shared string "run")
$ws.Range("A2").Value = "run"

# Update the active selection to H8 (matches sqref in diff)
$ws.Range("H8").Select()
